$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-type columns (B = date-as-text, C = id-as-text)
# to remain text rather than being auto-converted to a date serial / number.
$ws.Range("B7:C8").NumberFormat = "@"

$ws.Range("A7").Value = 1582761600
$ws.Range("B7").Value = "2020-02-27"
$ws.Range("C7").Value = "5295"
$ws.Range("D7").Value = "INNATURE"
$ws.Range("E7").Value = 0.6
$ws.Range("F7").Value = 0.6
$ws.Range("G7").Value = 0.5600000000000001
$ws.Range("H7").Value = 0.5600000000000001
$ws.Range("I7").Value = 2592500

$ws.Range("A8").Value = 1582848000
$ws.Range("B8").Value = "2020-02-28"
$ws.Range("C8").Value = "5295"
$ws.Range("D8").Value = "INNATURE"
$ws.Range("E8").Value = 0.55
$ws.Range("F8").Value = 0.57
$ws.Range("G8").Value = 0.52
$ws.Range("H8").Value = 0.545
$ws.Range("I8").Value = 4336600

# Drop the number-format override style again, so the new cells keep the
# default (unstyled) appearance, matching the rest of the data rows.
$ws.Range("B7:C8").Style = "Normal"
